# Change drivers for fugitive emissions sectors
# Rows 30-32 on the "Sectors" sheet are the fugitive-emissions sectors
# (1B1_Fugitive-solid-fuels, 1B2_Fugitive-petr-and-gas, 1B2d_Fugitive-other-energy).
# Their "activity" driver (column B) changes from population ("pop") to a
# fuel-supply based driver, and the "units" column (C) switches from a
# numeric placeholder (1000) to the text unit "kt".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# 1B2_Fugitive-petr-and-gas (row 31) and 1B2d_Fugitive-other-energy (row 32)
# both move to the "refinery-and-natural-gas" driver. Write these first so
# that string lands in the shared-string table before "coal-dom-supply".
$ws.Range("B31").ClearFormats()
$ws.Range("B31").Value = "refinery-and-natural-gas"
$ws.Range("C31").Value = "kt"

$ws.Range("B32").ClearFormats()
$ws.Range("B32").Value = "refinery-and-natural-gas"
$ws.Range("C32").Value = "kt"

# 1B1_Fugitive-solid-fuels (row 30) moves to the "coal-dom-supply" driver.
$ws.Range("B30").ClearFormats()
$ws.Range("B30").Value = "coal-dom-supply"
$ws.Range("C30").Value = "kt"

# Restore the single-cell selection on B30 (matches the saved view state).
$ws.Activate()
$ws.Range("B30").Select()
